# Applies the "Optuna Attempt (go back with original)" data refresh:
# updates MyForecast, Inventory Coverage, Stockout Risk and Seasonality Index
# values on the "Forecast Comparison" sheet, and the derived totals on the
# "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------

# Row 2 (W8)
$wsForecast.Range("D2").Value = 8
$wsForecast.Range("H2").Value = 0.49
$wsForecast.Range("I2").Value = "High"
$wsForecast.Range("L2").Value = 0.87

# Row 3 (W9)
$wsForecast.Range("D3").Value = 9
$wsForecast.Range("L3").Value = 1.1

# Row 4 (W10)
$wsForecast.Range("L4").Value = 1.06

# Row 5 (W11)
$wsForecast.Range("L5").Value = 0.85

# Row 7 (W13)
$wsForecast.Range("L7").Value = 1.02

# Row 8 (W14)
$wsForecast.Range("L8").Value = 0.89

# Row 9 (W15)
$wsForecast.Range("L9").Value = 0.96

# Row 10 (W16)
$wsForecast.Range("D10").Value = 8
$wsForecast.Range("L10").Value = 1.06

# Row 11 (W17)
$wsForecast.Range("D11").Value = 9
$wsForecast.Range("L11").Value = 0.92

# Row 12 (W18)
$wsForecast.Range("D12").Value = 7
$wsForecast.Range("L12").Value = 1

# Row 13 (W19)
$wsForecast.Range("D13").Value = 8
$wsForecast.Range("L13").Value = 0.97

# Row 14 (W20)
$wsForecast.Range("L14").Value = 0.99

# Row 15 (W21)
$wsForecast.Range("L15").Value = 0.9

# Row 16 (W22)
$wsForecast.Range("L16").Value = 1.08

# Row 17 (W23)
$wsForecast.Range("L17").Value = 1.16

# --- Summary sheet --------------------------------------------------------------
# These cells hold numeric-looking text (originally written as inline strings),
# so force literal text entry (leading apostrophe) and then restore the
# "Normal" style so no stray number-format is left behind on the cell.

$wsSummary.Range("B9").Value = "'147"
$wsSummary.Range("B9").Style = "Normal"

$wsSummary.Range("B10").Value = "'76"
$wsSummary.Range("B10").Style = "Normal"

$wsSummary.Range("B11").Value = "'37"
$wsSummary.Range("B11").Style = "Normal"

$wsSummary.Range("B14").Value = "'8"
$wsSummary.Range("B14").Style = "Normal"
